# Apply the "fixed workflow" re-windowing of the reaction-sensitivity
# tables on both sheets (NBR and BAR): slide the Cutoff window from
# [1..19] to [5..19] (still 15 rows) and refresh the Reaction_number
# counts, then drop the now-unused trailing rows.

$wb = $excel.ActiveWorkbook

# New B (Cutoff) and C (Reaction_number) values for rows 2..16.
$nbrB = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
$nbrC = @(130,131,133,133,128,126,126,126,126,127,127,126,127,127,126)

$barB = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
$barC = @(600,591,587,587,590,589,589,589,589,586,586,588,582,582,585)

foreach ($sheetIndex in 1..2) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    if ($sheetIndex -eq 1) {
        $bVals = $nbrB
        $cVals = $nbrC
    } else {
        $bVals = $barB
        $cVals = $barC
    }

    for ($i = 0; $i -lt 15; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $bVals[$i]
        $ws.Cells.Item($row, 3).Value = $cVals[$i]
    }

    # Rows 17-20 are no longer part of the window; remove them.
    $ws.Range("A17:A20").EntireRow.Delete()
}
